$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(1, 1).Value = "Datos actualizados a 12 de Junio de 2020 a las 21:15"
$ws.Cells.Item(4, 2).Value = 2107516
$ws.Cells.Item(4, 3).Value = 17815
$ws.Cells.Item(4, 4).Value = 820284
$ws.Cells.Item(4, 5).Value = 1170606
$ws.Cells.Item(4, 7).Value = 592
$ws.Cells.Item(4, 8).Value = 116626
$ws.Cells.Item(7, 2).Value = 309603
$ws.Cells.Item(7, 3).Value = 11320
$ws.Cells.Item(7, 4).Value = 154231
$ws.Cells.Item(7, 5).Value = 146482
$ws.Cells.Item(12, 2).Value = 187233
$ws.Cells.Item(12, 3).Value = 438
$ws.Cells.Item(12, 5).Value = 6772
$ws.Cells.Item(12, 7).Value = 10
$ws.Cells.Item(12, 8).Value = 8861
$ws.Cells.Item(78, 1).Value = "Costa de Marfil"
$ws.Cells.Item(78, 2).Value = 4684
$ws.Cells.Item(78, 3).Value = 280
$ws.Cells.Item(78, 4).Value = 2263
$ws.Cells.Item(78, 5).Value = 2376
$ws.Cells.Item(78, 7).Value = 4
$ws.Cells.Item(78, 8).Value = 45
$ws.Cells.Item(79, 1).Value = "Consejo Danes para los Refugiados"
$ws.Cells.Item(79, 2).Value = 4637
$ws.Cells.Item(79, 3).Value = 122
$ws.Cells.Item(79, 4).Value = 580
$ws.Cells.Item(79, 5).Value = 3956
$ws.Cells.Item(79, 7).Value = 3
$ws.Cells.Item(79, 8).Value = 101
$ws.Cells.Item(80, 1).Value = "Republica de Yibuti"
$ws.Cells.Item(80, 2).Value = 4441
$ws.Cells.Item(80, 3).Value = 43
$ws.Cells.Item(80, 4).Value = 2730
$ws.Cells.Item(80, 5).Value = 1673
$ws.Cells.Item(80, 7).Value = 1
$ws.Cells.Item(80, 8).Value = 38
$ws.Cells.Item(81, 2).Value = 4426
$ws.Cells.Item(81, 3).Value = 54
$ws.Cells.Item(81, 4).Value = 3106
$ws.Cells.Item(81, 5).Value = 1296
$ws.Cells.Item(81, 7).Value = 1
$ws.Cells.Item(81, 8).Value = 24
$ws.Cells.Item(108, 1).Value = "Costa Rica"
$ws.Cells.Item(108, 2).Value = 1612
$ws.Cells.Item(108, 3).Value = 74
$ws.Cells.Item(108, 4).Value = 731
$ws.Cells.Item(108, 5).Value = 869
$ws.Cells.Item(108, 8).Value = 12
$ws.Cells.Item(109, 1).Value = "Mauritania"
$ws.Cells.Item(109, 2).Value = 1572
$ws.Cells.Item(109, 3).Value = 133
$ws.Cells.Item(109, 4).Value = 278
$ws.Cells.Item(109, 5).Value = 1213
$ws.Cells.Item(109, 7).Value = 7
$ws.Cells.Item(109, 8).Value = 81
$ws.Cells.Item(110, 1).Value = "Eslovaquia"
$ws.Cells.Item(110, 2).Value = 1542
$ws.Cells.Item(110, 3).Value = 1
$ws.Cells.Item(110, 4).Value = 1409
$ws.Cells.Item(110, 5).Value = 105
$ws.Cells.Item(110, 8).Value = 28
$ws.Cells.Item(111, 1).Value = "Nueva Zelanda"
$ws.Cells.Item(111, 2).Value = 1504
$ws.Cells.Item(111, 3).Value = 0
$ws.Cells.Item(111, 4).Value = 1482
$ws.Cells.Item(111, 5).Value = 0
$ws.Cells.Item(111, 8).Value = 22
$ws.Cells.Item(112, 1).Value = "Eslovenia"
$ws.Cells.Item(112, 2).Value = 1490
$ws.Cells.Item(112, 3).Value = 2
$ws.Cells.Item(112, 4).Value = 1359
$ws.Cells.Item(112, 5).Value = 22
$ws.Cells.Item(112, 8).Value = 109
$ws.Cells.Item(113, 1).Value = "Nicaragua"
$ws.Cells.Item(113, 2).Value = 1464
$ws.Cells.Item(113, 4).Value = 953
$ws.Cells.Item(113, 5).Value = 456
$ws.Cells.Item(113, 8).Value = 55
$ws.Cells.Item(127, 2).Value = 978
$ws.Cells.Item(127, 3).Value = 4
$ws.Cells.Item(127, 4).Value = 881
$ws.Cells.Item(127, 5).Value = 32
$ws.Cells.Item(206, 1).Value = "Islas Malvinas"
$ws.Cells.Item(207, 1).Value = "Groenlandia"
$ws.Cells.Item(208, 1).Value = "Islas Turcas y Caicos"
$ws.Cells.Item(208, 4).Value = 11
$ws.Cells.Item(208, 8).Value = 1
$ws.Cells.Item(209, 1).Value = "Santa Sede"
$ws.Cells.Item(209, 4).Value = 12
$ws.Cells.Item(209, 8).Value = 0
$ws.Cells.Item(210, 1).Value = "Seychelles"
$ws.Cells.Item(210, 4).Value = 11
$ws.Cells.Item(210, 8).Value = 0
$ws.Cells.Item(211, 1).Value = "Montserrat"
$ws.Cells.Item(211, 4).Value = 10
$ws.Cells.Item(211, 8).Value = 1
$ws.Cells.Item(213, 1).Value = "Papua Nueva Guinea"
$ws.Cells.Item(213, 4).Value = 8
$ws.Cells.Item(213, 8).Value = 0
$ws.Cells.Item(214, 1).Value = "Islas Virgenes Britanicas"
$ws.Cells.Item(214, 4).Value = 7
$ws.Cells.Item(214, 8).Value = 1
